$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header cells (e.g. H1 - bold, bordered, centered).
$ws.Range("I1").Value = "I0"
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)

$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)

# Row 2 is a special case: I2 = 7, J2 = 7
$ws.Cells.Item(2, 9).Value = 7
$ws.Cells.Item(2, 10).Value = 7

# Rows 3-39: I = 1, J = value currently in column H (IP)
for ($r = 3; $r -le 39; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
